$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the two labels (trailing space is intentional, matches target data)
$ws.Range("A5").Value = "Beitragsbemessungsgrenze GKV "
$ws.Range("A6").Value = "Jahresarbeitsentgeltgrenze GKV "

# Update the selected cell shown in the sheet view
$ws.Activate()
$ws.Range("A4").Select()
